function Clean-Text($s) {
    # Replace any HTML tag (e.g. <p>, </p>, <strong>, </strong>, <br>) with a single space
    $t = [System.Text.RegularExpressions.Regex]::Replace($s, "<[^>]+>", " ")
    # Collapse any run of whitespace (spaces, double-spaces, newlines, tabs) into one space
    $t = [System.Text.RegularExpressions.Regex]::Replace($t, "\s+", " ")
    # Trim leading/trailing whitespace
    $t = $t.Trim()
    return $t
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$firstRow = $used.Row
$firstCol = $used.Column
$rowCount = $used.Rows.Count
$colCount = $used.Columns.Count
$lastRow = $firstRow + $rowCount - 1
$lastCol = $firstCol + $colCount - 1

for ($r = $firstRow; $r -le $lastRow; $r++) {
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $v = $cell.Value2
        if ($v -is [string]) {
            if ($v.Length -gt 0) {
                $cleaned = Clean-Text $v
                if ($cleaned -ne $v) {
                    $cell.Value = $cleaned
                }
            }
        }
    }
}
